$d = $word.ActiveDocument

# Locate the title run that currently reads "progressãoUFBA".
$found = $d.Content
$ok = $found.Find.Execute("progressãoUFBA", $true, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0)

if (-not $ok) {
    throw "Could not find 'progressãoUFBA' in the document."
}

$target = $d.Range($found.Start, $found.End)

# Rebuild the same text ("ProgressãoUFBA": first letter capitalised) but as two
# runs with identical run properties, exactly like the authored edit, by
# replacing the range's contents via WordprocessingML rather than plain text
# (plain .Text assignment gets coalesced back into a single run).
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p>' +
       '<w:r><w:rPr><w:sz w:val="72"/><w:szCs w:val="72"/></w:rPr><w:t>P</w:t></w:r>' +
       '<w:r><w:rPr><w:sz w:val="72"/><w:szCs w:val="72"/></w:rPr><w:t>rogressãoUFBA</w:t></w:r>' +
       '</w:p>' +
       '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData>' +
       '</pkg:part>' +
       '</pkg:package>'

$target.InsertXML($xml)
